$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

$newValues = @(
    "67-27=",
    "44-6=",
    "41-23=",
    "22+76=",
    "88-75=",
    "64+27=",
    "83-11=",
    "34+38=",
    "68-38=",
    "55-54=",
    "53+37=",
    "34-19=",
    "10+73=",
    "35+8=",
    "87-15=",
    "43+26=",
    "86-5=",
    "29+37=",
    "47+35=",
    "86+7=",
    "51-41=",
    "1+92=",
    "47+32=",
    "25+64=",
    "29+66=",
    "85-26=",
    "63-34=",
    "63+9=",
    "55+14=",
    "21+24=",
    "4+69=",
    "54+14=",
    "74-13=",
    "80-7=",
    "54-23=",
    "11-7=",
    "52-28=",
    "15-0=",
    "58-25=",
    "73-65=",
    "39-27=",
    "59-47=",
    "45+27=",
    "35+25=",
    "75-43=",
    "31+66=",
    "4+64=",
    "19+7=",
    "38+35=",
    "69+17=",
    "42+23=",
    "85+8=",
    "95+2=",
    "91-45=",
    "64-44=",
    "87-62=",
    "73-69=",
    "46-0=",
    "5+16=",
    "90-47=",
    "79-42=",
    "82-16=",
    "60-32=",
    "24+11=",
    "38+6=",
    "53-22=",
    "17+74=",
    "68-23=",
    "83-38=",
    "50+18=",
    "41+35=",
    "26+14=",
    "83-31=",
    "65-9=",
    "30+18=",
    "59-7=",
    "12+5=",
    "25+25=",
    "86-71=",
    "8+19=",
    "90-40=",
    "9+10=",
    "71+12=",
    "54-19=",
    "91-65=",
    "6+78=",
    "79-36=",
    "82+6=",
    "43-35=",
    "36+60=",
    "16+56=",
    "16+0=",
    "16+3=",
    "16+44=",
    "92-66=",
    "48-38=",
    "92-86=",
    "15+56=",
    "32+1=",
    "33+26="
)

$cols = 5
$idx = 0
for ($r = 1; $r -le $tbl.Rows.Count; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $tbl.Cell($r, $c)
        $cell.Range.Text = $newValues[$idx]
        $idx++
    }
}

Write-Output "Updated $idx cells"